# Insert a new row at position 42 (shifts existing rows 42-73 down to 43-74)
# and populate it with a new Pepino dulce price record, matching the
# author's weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 45072
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = 100112043
$ws.Cells.Item(42, 7).Value = "Pepino dulce"
$ws.Cells.Item(42, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 290
$ws.Cells.Item(42, 11).Value = 13000
$ws.Cells.Item(42, 12).Value = 14000
$ws.Cells.Item(42, 13).Value = 13517
$ws.Cells.Item(42, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 751
$ws.Cells.Item(42, 17).Value = 18
$ws.Cells.Item(42, 18).Value = "Hortaliza"
